$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SBER")

# Shift every price in the grid (A2:A42) up by 7, preserving the
# original floating point noise already baked into the cells.
$rng = $ws.Range("A2:A42")
foreach ($cell in $rng.Cells) {
    $cell.Value2 = $cell.Value2 + 7
}

# Select the full price column and make A2 the active cell, matching
# the new selection left behind after the bot edit.
$rng.Select()
$excel.ActiveCell = $ws.Range("A2")
